$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (column F) values for three events
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 1194
$ws1.Range("F16").Value = 129
$ws1.Range("F18").Value = 178

# Sheet "全部类型": same three events appear one row lower
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 1194
$ws4.Range("F17").Value = 129
$ws4.Range("F19").Value = 178
